$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 852; this shifts the existing rows 852-908
# down to 853-909 (matching the rest of the diff, which is just every
# subsequent row's content moving down by one row).
$ws.Rows(852).Insert()

# Populate the newly inserted row 852 with its data.
$ws.Range("A852").Value = 10
$ws.Range("B852").Value = "Vega Modelo de Temuco"
$ws.Range("C852").Value = "La Araucanía"
$ws.Range("D852").Value = 44826
$ws.Range("E852").Value = 9
$ws.Range("F852").Value = 100112021
$ws.Range("G852").Value = "Ají"
$ws.Range("H852").Value = "Inferno"
$ws.Range("I852").Value = "Primera"
$ws.Range("J852").Value = 80
$ws.Range("K852").Value = 17000
$ws.Range("L852").Value = 17000
$ws.Range("M852").Value = 17000
$ws.Range("N852").Value = '$/caja 15 kilos'
$ws.Range("O852").Value = "Región de Arica y Parinacota"
$ws.Range("P852").Value = 1133
$ws.Range("Q852").Value = 15
$ws.Range("R852").Value = "Hortaliza"
